$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "aa"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "tepong32@gmail.com"

$ws.Range("A5").Value = "aa         a"
$ws.Range("B5").Value = 33
$ws.Range("C5").Value = "tepong32@gmail.com"
